# Added area column to identify printers:
# abbreviate a few long "Area" values in column D and add a
# "Siglas / Significado" (Acronym / Meaning) legend table in columns F:G
# explaining those abbreviations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Abbreviate "Creditos y cobranzas" and build its legend entry ---
$ws.Range("D7").Value = "Cred. Cob."

$ws.Range("F1").Value = "Siglas"
$ws.Range("G1").Value = "Significado"

$ws.Range("F2").Value = "Cred. Cob."
$ws.Range("G2").Value = "Creditos y Cobranzas"

# --- Abbreviate "Gestion Humana" and build its legend entry ---
$ws.Range("D13").Value = "G.H."

$ws.Range("F3").Value = "G.H."
$ws.Range("G3").Value = "Gestion Humana"

# --- Abbreviate "San Juan de Yapacany" and build its legend entry ---
$ws.Range("D21").Value = "S.J. Yapacany"

$ws.Range("F4").Value = "S.J. Yapacany"
$ws.Range("G4").Value = "San Juan de Yapacany"

# Apply the same header style as A1:D1 to F1:G1
$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122) # xlPasteFormats

# Apply the same body style as A2:D4 to F2:G4
$ws.Range("A2").Copy()
$ws.Range("F2:G4").PasteSpecial(-4122) # xlPasteFormats

# Row heights for rows 2-4 (wrap text legend rows)
$ws.Rows("2:4").RowHeight = 25.5

# Update the active selection to match final state
$ws.Range("G5").Select()
